$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New commit row (row 30): commit message text in C30, hours spent in G30
$ws.Range("C30").Value = "Runner & Merchant completly working"
$ws.Range("G30").Value = 4

# Total(h) now needs to sum through the new row
$ws.Range("G39").Formula = "=SUM(G4:G30)"

# Match the saved view/selection state (scrolled down, C31 selected)
$ws.Range("C31").Select()
$excel.ActiveWindow.ScrollRow = 13
